$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.498.04'
$ws.Range('E2').Value = '  +1.36%  '

$ws.Range('D3').Value = '3.704.86'
$ws.Range('E3').Value = '  +0.72%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '617.26'
$ws.Range('E5').Value = '  +6.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '196.58'
$ws.Range('E6').Value = '  +14.76%  '

$ws.Range('E7').Value = '  +2.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.726'
$ws.Range('E9').Value = '  +3.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '60.50'
$ws.Range('E10').Value = '  +17.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.160'
$ws.Range('E11').Value = '  -0.47%  '

$ws.Range('E12').Value = '  -1.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '10.44'
$ws.Range('E13').Value = '  +0.21%  '

$ws.Range('D14').Value = '4.298.27'
$ws.Range('E14').Value = '  +0.66%  '

$ws.Range('D15').Value = '3.703.64'
$ws.Range('E15').Value = '  +0.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '19.41'
$ws.Range('E16').Value = '  +0.59%  '

$ws.Range('E17').Value = '  +0.83%  '

$ws.Range('E18').Value = '  +2.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '12.86'
$ws.Range('E19').Value = '  +0.33%  '

$ws.Range('D20').Value = '68.383.06'
$ws.Range('E20').Value = '  +1.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '408.64'
$ws.Range('E21').Value = '  +0.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '4.66'
$ws.Range('E22').Value = '  +4.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '89.92'
$ws.Range('E23').Value = '  +2.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '3.07'
$ws.Range('E24').Value = '  +1.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '11.47'
$ws.Range('E25').Value = '  +8.36%  '

$ws.Range('E26').Value = '  +2.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '6.03'
$ws.Range('E27').Value = '  +0.85%  '

$ws.Range('E28').Value = '  +2.00%  '

$ws.Range('E29').Value = '  +2.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '32.69'
$ws.Range('E30').Value = '  +0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '7.62'
$ws.Range('E31').Value = '  +2.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '12.68'
$ws.Range('E32').Value = '  +2.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '46.51'
$ws.Range('E33').Value = '  +8.24%  '

$ws.Range('E34').Value = '  +5.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '634.23'
$ws.Range('E35').Value = '  +6.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '67.64'
$ws.Range('E36').Value = '  +4.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.412'
$ws.Range('E37').Value = '  +4.30%  '

$ws.Range('D38').Value = '0.0₃0821'
$ws.Range('E38').Value = '  -7.00%  '

$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('E40').Value = '  -0.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '0.139'
$ws.Range('E41').Value = '  +4.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '3.04'
$ws.Range('E42').Value = '  +2.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '0.0444'
$ws.Range('E43').Value = '  +2.10%  '

$ws.Range('E44').Value = '  -1.19%  '

$ws.Range('D45').Value = '2.915.59'
$ws.Range('E45').Value = '  +4.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '9.41'
$ws.Range('E46').Value = '  +2.90%  '

$ws.Range('E47').Value = '  +4.87%  '

$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '146.87'
$ws.Range('E48').Value = '  +2.72%  '

$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '2.70'
$ws.Range('E49').Value = '  +0.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '3.08'
$ws.Range('E50').Value = '  -2.44%  '

$ws.Range('E51').Value = '  -12.23%  '
